# Split the single long <w:t> runs in the "Programa" and "Bibliografia"
# paragraphs into multiple <w:t> elements separated by manual line breaks
# (<w:br/>), one per numbered item, matching the author's re-formatting.

$d = $word.ActiveDocument

function Insert-LineBreak {
    param(
        [string]$FindText,
        [string]$ReplaceText
    )
    $d.Content.Find.Execute($FindText, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $ReplaceText, 2) | Out-Null
}

# --- "Programa" paragraph -------------------------------------------------
Insert-LineBreak "e Aços. Demais Processos.2) Process" "e Aços. Demais Processos.^l2) Process"
Insert-LineBreak "ndição. Demais Processos.3) Process" "ndição. Demais Processos.^l3) Process"
Insert-LineBreak "sinagem. Demais Processos4) Process" "sinagem. Demais Processos^l4) Process"

# --- "Bibliografia" paragraph ----------------------------------------------
Insert-LineBreak "ABM, São Paulo, 20072. Fathi" "ABM, São Paulo, 2007^l2. Fathi"
Insert-LineBreak "e Publishers, 1986. 3. Luiz " "e Publishers, 1986. ^l3. Luiz "
Insert-LineBreak "a, São Paulo, 1997. 4. Alan " "a, São Paulo, 1997. ^l4. Alan "
Insert-LineBreak "nkian, Lisboa, 1975.5. ASM H" "nkian, Lisboa, 1975.^l5. ASM H"
Insert-LineBreak "gy P.R. Beeley, 19726. John " "gy P.R. Beeley, 1972^l6. John "
Insert-LineBreak "orth-Heinemann, 19917. M. Si" "orth-Heinemann, 1991^l7. M. Si"
Insert-LineBreak "ABM, S.Paulo, 1979. 8. Amaur" "ABM, S.Paulo, 1979. ^l8. Amaur"
Insert-LineBreak ", Campinas, SP, 20089. Mauri" ", Campinas, SP, 2008^l9. Mauri"
Insert-LineBreak "cos, Rio de Janeiro.10. AVIT" "cos, Rio de Janeiro.^l10. AVIT"
Insert-LineBreak "ed; New Delhi, 1977.11. RODR" "ed; New Delhi, 1977.^l11. RODR"
Insert-LineBreak "ra, v.1 e v.2, 2010.12. CETL" "ra, v.1 e v.2, 2010.^l12. CETL"
Insert-LineBreak "r Ltda, 260p., 2005.13. BRES" "r Ltda, 260p., 2005.^l13. BRES"
Insert-LineBreak "mp, v.1 e v.2, 1986.14. DINI" "mp, v.1 e v.2, 1986.^l14. DINI"
